$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the August data row (row 6)
$ws.Range("A6").Value = 45169
$ws.Range("A6").Style = $ws.Range("A2").Style
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B6").Value = 2143000
